# test_excel_ui.xlsx: "Changed column order in Instruments tab"
#
# The "Instruments" worksheet had its "Time Channel" (column E) and
# "Fluorescence Channels" (column F) columns swapped, and the Instruments
# sheet became the active/selected tab (it previously was "Samples").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instruments")

# Swap columns E and F as whole columns, so values, styles and the
# per-column width travel together exactly like a drag-to-reorder in the
# Excel UI: cut column F ("Fluorescence Channels") and insert it in front
# of column E ("Time Channel"), which pushes the old E into the F slot.
$ws.Columns.Item(6).Cut()
$ws.Columns.Item(5).Insert()

# The Instruments tab is now the one the user is looking at.
$ws.Activate()
